# Update Test Suite Statistics sheet for SourceProviderRDO - new test case written.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 ("SourceProviderRDO") total test case counts bump from 8 to 9.
$ws.Range("B3").Value = 9
$ws.Range("C3").Value = 9

# Recalculate dependent SUM formulas (G4/G5) to reflect the new totals.
$excel.Calculate()

# Move the active selection to D3, matching the saved view state.
$ws.Range("D3").Select()

$wb.Save()
